# Apply edits to LOM3018 worksheet per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Variables holding the new / changed text values ---
$ativacaoDate = '01/01/2022'
$objetivosText = 'A disciplina busca introduzir o aluno ao ambiente de engenharia, propondo problemas desafiadores gerando aptidão para solução de problemas. Apresentar a Engenharia de Materiais e seus campos de atuação, aspectos legais e éticos, bem como o mercado de trabalho para o engenheiro de materiais no Século XXI. Propiciar aos alunos uma visão geral do curso, com apresentação do currículo do curso de Engenharia de Materiais da EEL. Apresentar aos alunos uma visão da evolução histórica dos materiais com o homem. Descrever exemplos marcantes da introdução de novos materiais e as mudanças sociais provocadas. Apresentar o caráter interdisciplinar da Ciência e Engenharia de Materiais e suas ligações com outros ramos da Ciência. Apresentar estudos de caso demonstrando este caráter interdisciplinar.'
$newTeacher = '7459752 - Maria Ismenia Sodero Toledo Faria'
$programaResumidoText = '1- A importância dos materiais na evolução do homem na pré-história. Alquimia, Revolução Científica e a Revolução Industrial. 2-O Engenheiro como um profissional, funções da engenharia, a ética e comunicação na engenharia 3-A grandes áreas da Engenharia de Materiais. A interdisciplinaridade da Ciência e Engenharia de Materiais. 4- Perspectivas para a Engenharia de Materiais no século XXI. 5- O currículo do curso de engenharia de materiais da EEL-USP. 6- Noções básicas de Projetos em Engenharia.Em todos o conteúdo do curso serão abordados aspectos sociais, ambientais, éticos, legais e econômicos para ampliar as competências dos alunos'
$programaText = '1- As características importantes de um engenheiro: aptidões interpessoais, aptidões de comunicação, liderança e competência. O engenheiro, profissional que busca solucionar problemas. 2-A Engenharia de Materiais: áreas de atuação e mercado de trabalho. Aplicação. A importância dos materiais na evolução do homem, as grandes áreas e interdisciplinaridade da Ciência e Engenharia de Materiais. Visita ao Departamento de Engenharia de Materiais. Conhecimento dos Grupos de Pesquisa do Departamento. Perspectivas para a Engenharia de Materiais no século XXI. 3- O campo de trabalho do engenheiro de materiais e suas áreas de atuação. Visita externa para integralização dos conhecimentos. 4- O currículo do curso de engenharia de materiais na EEL/USP. 5- Apresentação do método de trabalho com projetos, definindo os atributos de um projeto de engenharia, mapas conceituais e ferramentas que ilustram ideias e relações entre elas. Formular estratégias para resolução de problemas de engenharia. Estudo de casos'
$metodoText = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras'
$criterioText = 'Média Aritmética dos Projetos, Trabalhos, Relatórios e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.'
$normaRecupText = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$bibliografiaText = '1) BROCKMAN, J.B. Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2010.2) M.T. HOLTZAPPLE, W.D. REECE, Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2006.2) CALLISTER Jr., W.D. Ciência e Engenharia de Materiais: Uma Introdução. LTC Livros Científicos Editora, 7a.ed., 2008. 4) - COHEN, M. (Ed.). Ciência e Engenharia de Materiais: sua Evolução, Prática e Perspectivas. Parte I: Materiais na história e na sociedade, 98p. Parte II: A Ciência e Engenharia de Materiais como uma multidisciplina, Tradução: José Roberto Gonçalves da Silva, São Carlos, UFSCar, 1985.5) Artigos científicos'

# --- 1. Update 'Ativacao:' date (row 8, B and C) keeping it as plain text ---
# Use a formula-then-paste-values trick so Excel doesn't auto-convert the
# dd/mm/yyyy-looking string into a date serial number.
$ws.Range("B8").Formula = '="' + $ativacaoDate + '"'
$ws.Range("C8").Formula = '="' + $ativacaoDate + '"'
$ws.Range("B8:C8").Copy()
$ws.Range("B8:C8").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- 2. Update 'Objetivos:' body text (row 10, B and C) ---
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

# --- 3. Insert a new row for the second 'Docentes responsaveis' entry ---
$ws.Rows("14").Insert()
$ws.Range("B14").Value = $newTeacher
$ws.Range("C14").Value = $newTeacher

# --- 4. Update 'Programa resumido:' body text (now row 15, B and C) ---
$ws.Range("B15").Value = $programaResumidoText
$ws.Range("C15").Value = $programaResumidoText

# --- 5. Update 'Programa:' body text (now row 17, B and C) ---
$ws.Range("B17").Value = $programaText
$ws.Range("C17").Value = $programaText

# --- 6. Update 'Metodo:' body text (now row 20, B and C) ---
$ws.Range("B20").Value = $metodoText
$ws.Range("C20").Value = $metodoText

# --- 7. Update 'Criterio:' body text (now row 21, B and C) ---
$ws.Range("B21").Value = $criterioText
$ws.Range("C21").Value = $criterioText

# --- 8. Update 'Norma de recuperacao:' body text (now row 22, B and C) ---
$ws.Range("B22").Value = $normaRecupText
$ws.Range("C22").Value = $normaRecupText

# --- 9. Update 'Bibliografia:' body text (now row 23, B and C) ---
$ws.Range("B23").Value = $bibliografiaText
$ws.Range("C23").Value = $bibliografiaText

